## Fruta / hortaliza, semanal
## Insert one new weekly price-report row above the current row 204
## (Feria Lagunitas de Puerto Montt - Uva sheet). Excel shifts every row
## from 204..266 down to 205..267 and grows the used range to A1:T267.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 204 - everything below shifts down one row.
$ws.Rows(204).Insert()

# Populate the newly inserted row 204 with the new record.
$ws.Cells.Item(204, 1).Value2  = 4
$ws.Cells.Item(204, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(204, 3).Value2  = "Los Lagos"
$ws.Cells.Item(204, 4).Value2  = 44917
$ws.Cells.Item(204, 5).Value2  = 10
$ws.Cells.Item(204, 6).Value2  = "Fruta"
$ws.Cells.Item(204, 7).Value2  = 100109
$ws.Cells.Item(204, 8).Value2  = "Uva"
$ws.Cells.Item(204, 9).Value2  = 100109001
$ws.Cells.Item(204, 10).Value2 = "Uva"
$ws.Cells.Item(204, 11).Value2 = "Flame Seedless"
$ws.Cells.Item(204, 12).Value2 = "Primera"
$ws.Cells.Item(204, 13).Value2 = 200
$ws.Cells.Item(204, 14).Value2 = 20000
$ws.Cells.Item(204, 15).Value2 = 22000
$ws.Cells.Item(204, 16).Value2 = 21000
$ws.Cells.Item(204, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(204, 18).Value2 = "Provincia de Copiapó"
$ws.Cells.Item(204, 19).Value2 = 2100
$ws.Cells.Item(204, 20).Value2 = 10
